# Zeitaufzeichnung update: add a new day of work (row 31) and correct the
# hours recorded for the previous entry (row 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: hours worked corrected from 3.5 to 4.5
$ws.Range("A30").Value = 4.5

# Row 31: new entry - 4.5 hours on 2020-01-15, 17:00-21:30, extra presentation work
# Copy the date formatting from the row above so B31 keeps the same date
# number format (rather than Excel inventing a brand new date/time format).
$ws.Range("B30").Copy()
$ws.Range("B31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A31").Value = 4.5
$ws.Range("B31").Value = (Get-Date -Year 2020 -Month 1 -Day 15 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C31").Value = "17:00-21:30"
$ws.Range("D31").Value = "Präs. zusatz, best practice"

# Active selection moves to A31
$ws.Range("A31").Select()
